# "functionifying the script to facilitate reading"
#
# On the "inputs" sheet, add a "Watershed" column (between "Species" and
# "Waterbody") and two new example rows describing watershed-level /
# multi-species model runs.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("inputs")

# Insert a new column at C ("Watershed"), pushing the existing "Waterbody"
# column to D and "Established_in_Waterbody" to E.
$ws1.Columns.Item(3).Insert()
$ws1.Range("C1").Value2 = "Watershed"

# New scenario row: run the model across a list of species for every
# waterbody in the Columbia River Watershed.
$ws1.Range("C15").Value2 = "Columbia River Watershed"
$ws1.Range("D15").Value2 = "NA"
$ws1.Range("E15").Value2 = "NA"
$ws1.Range("B15").Value2 = "Smallmouth bass, Goldfish, Pumpkinseed"
$ws1.Range("F15").Value2 = "This would run the model for each species in the list for all waterbodies in the chosen watershed THAT HAVE OVERLAPPED WITH SAR (preliminary step); priority ranking is done ACROSS species. Also, have priority columns for MONITORING, DOWNSTREAM EXTENT"
$ws1.Range("F15").WrapText = $true
$ws1.Rows.Item(15).RowHeight = 57.6

# Same idea, but for the Fraser River Watershed.
$ws1.Range("B16").Value2 = "Smallmouth bass, Goldfish, Pumpkinseed"
$ws1.Range("C16").Value2 = "Fraser River Watershed"
$ws1.Range("D16").Value2 = "NA"
$ws1.Range("E16").Value2 = "NA"

# Widen columns so the new/longer content fits.
$ws1.Columns.Item(1).ColumnWidth = 19.833333333333332
$ws1.Columns.Item(2).ColumnWidth = 36.916666666666664
$ws1.Columns.Item(3).ColumnWidth = 23.5
$ws1.Columns.Item(5).ColumnWidth = 16.25
$ws1.Columns.Item(6).ColumnWidth = 62.5

# Leave the "species_predvars" sheet's own selection updated too (visited
# while reviewing the new watershed scenarios), then return focus to the
# "inputs" sheet where the edits were made.
$ws4 = $wb.Worksheets.Item("species_predvars")
$ws4.Activate()
[void]$ws4.Range("K17").Select()

$ws1.Activate()
[void]$ws1.Range("C17").Select()
